$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Column A holds numeric-looking IDs that must stay stored as *text*
# (matching every existing row). Writing "22" straight into a General-
# formatted cell makes Excel infer a number, so: format as Text, assign
# the value, then clear the formatting override back to the sheet's
# default (General) so no stray style is left behind on the cell.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "22"
$ws.Range("A18").ClearFormats()

$ws.Range("B18").Value = "jhasgcdahjsgdhajsgdhas"
$ws.Range("C18").Value = "open"
$ws.Range("D18").Value = "2025-03-26T06:42:33Z"
$ws.Range("E18").Value = "bug"
